$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column widths to match the new layout (A: bestFit country-name column,
# B: default, C: wider for descriptive text)
$ws.Columns.Item(1).ColumnWidth = 43.67
$ws.Columns.Item(2).ColumnWidth = 8.33
$ws.Columns.Item(3).ColumnWidth = 24.1

# Append the new continent / region summary rows (254-275) that were
# pasted into the sheet. Column A is highlighted with a yellow fill to
# flag the newly added data; only columns A and E are populated.
$ws.Range("A254").Value = 'Central America'
$ws.Range("A254").Interior.Color = 65535
$ws.Range("E254").Value = 'North America'

$ws.Range("A255").Value = 'China Hong Kong SAR'
$ws.Range("A255").Interior.Color = 65535
$ws.Range("E255").Value = 'Asia'

$ws.Range("A256").Value = 'Eastern Africa'
$ws.Range("A256").Interior.Color = 65535
$ws.Range("E256").Value = 'Africa'

$ws.Range("A257").Value = 'Middle Africa'
$ws.Range("A257").Interior.Color = 65535
$ws.Range("E257").Value = 'Africa'

$ws.Range("A258").Value = 'Other Asia Pacific'
$ws.Range("A258").Interior.Color = 65535
$ws.Range("E258").Value = 'Asia'

$ws.Range("A259").Value = 'Other Caribbean'
$ws.Range("A259").Interior.Color = 65535
$ws.Range("E259").Value = 'South America'

$ws.Range("A260").Value = 'Other CIS'
$ws.Range("A260").Interior.Color = 65535
$ws.Range("E260").Value = 'Asia'

$ws.Range("A261").Value = 'Other Europe'
$ws.Range("A261").Interior.Color = 65535
$ws.Range("E261").Value = 'Europe'

$ws.Range("A262").Value = 'Other Middle East'
$ws.Range("A262").Interior.Color = 65535
$ws.Range("E262").Value = 'Asia'

$ws.Range("A263").Value = 'Other Northern Africa'
$ws.Range("A263").Interior.Color = 65535
$ws.Range("E263").Value = 'Africa'

$ws.Range("A264").Value = 'Other South America'
$ws.Range("A264").Interior.Color = 65535
$ws.Range("E264").Value = 'Africa'

$ws.Range("A265").Value = 'Other Southern Africa'
$ws.Range("A265").Interior.Color = 65535
$ws.Range("E265").Value = 'Africa'

$ws.Range("A266").Value = 'Russian Federation'
$ws.Range("A266").Interior.Color = 65535
$ws.Range("E266").Value = 'Asia'

$ws.Range("A267").Value = 'Total Africa'
$ws.Range("A267").Interior.Color = 65535
$ws.Range("E267").Value = 'Africa'

$ws.Range("A268").Value = 'Total Asia Pacific'
$ws.Range("A268").Interior.Color = 65535
$ws.Range("E268").Value = 'Asia'

$ws.Range("A269").Value = 'Total CIS'
$ws.Range("A269").Interior.Color = 65535
$ws.Range("E269").Value = 'Asia'

$ws.Range("A270").Value = 'Total Europe'
$ws.Range("A270").Interior.Color = 65535
$ws.Range("E270").Value = 'Europe'

$ws.Range("A271").Value = 'Total Middle East'
$ws.Range("A271").Interior.Color = 65535
$ws.Range("E271").Value = 'Asia'

$ws.Range("A272").Value = 'Total North America'
$ws.Range("A272").Interior.Color = 65535
$ws.Range("E272").Value = 'North America'

$ws.Range("A273").Value = 'Total S. & Cent. America'
$ws.Range("A273").Interior.Color = 65535
$ws.Range("E273").Value = 'South America'

$ws.Range("A274").Value = 'Total World'
$ws.Range("A274").Interior.Color = 65535
$ws.Range("E274").Value = 'Total World'
$ws.Range("E274").Interior.Color = 65535

$ws.Range("A275").Value = 'Western Africa'
$ws.Range("A275").Interior.Color = 65535
$ws.Range("E275").Value = 'Africa'
